# Add a new data row (row 3) to the download report and resize the
# first two columns to fit the new, wider content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New report line: date, resource name, download count
$ws.Range("A3").Value = "1.8.2023 0:00:00"
$ws.Range("B3").Value = "C# 10 in a Nutshell The Definitive Reference Joseph Albahari"
$ws.Range("C3").Value = 1

# Widen columns A and B so the new values are fully visible.
$ws.Columns(1).ColumnWidth = 14.627911249796549
$ws.Columns(2).ColumnWidth = 52.598978678385414
